$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"9.318123435519965e-06"
$ws.Range("C2").Value = [double]"0.0001537489499301437"
$ws.Range("D2").Value = [double]"157.8057217802531"
$ws.Range("E2").Value = [double]"246.9852506941017"
$ws.Range("G2").Value = [double]"404.7911355414282"

$ws.Range("B3").Value = [double]"1.505614041169197"
$ws.Range("C3").Value = [double]"86.29678392075563"
$ws.Range("D3").Value = [double]"3.082599426703578"
$ws.Range("E3").Value = [double]"6.48142807727062"
$ws.Range("G3").Value = [double]"97.36642546589903"

$ws.Range("B4").Value = [double]"0.06328177979961902"
$ws.Range("C4").Value = [double]"87981.0709163148"
$ws.Range("D4").Value = [double]"116886.6739907443"
$ws.Range("E4").Value = [double]"198602002.3250627"
$ws.Range("G4").Value = [double]"198806870.1332515"
